$d = $word.ActiveDocument

# 1) Insert a brand-new, truly-empty paragraph before the existing one.
#    (Paragraphs(1).Range.InsertParagraphBefore() would leave a stray
#    empty <w:r/> behind, so we splice in raw OOXML for a clean <w:p/>.)
$pkgEmptyPara = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p/></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$startRng = $d.Range(0, 0)
[void]$startRng.InsertXML($pkgEmptyPara)

# 2) In the (now second) paragraph, type the new leading text "sfjjsafljla"
#    ahead of the existing "sljlfs;fl" text.
$newText = "sfjjsafljla"
$para = $d.Paragraphs(2)
$rng = $para.Range
$rng.Collapse(1)
$insertStart = $rng.Start
$rng.InsertBefore($newText)
$insertEnd = $insertStart + $newText.Length

# 3) Re-home the _GoBack bookmark so it sits between the two runs
#    (right after the newly-typed text), matching the target diff,
#    instead of where it drifted to (the very end of the paragraph).
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$newBmRng = $d.Range($insertEnd, $insertEnd)
$d.Bookmarks.Add("_GoBack", $newBmRng)
